$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = '281474991205262-1743276209056'
$ws.Cells.Item(2, 2).Value = 'Mobile Usage'
$ws.Cells.Item(2, 3).Value = '2025-03-29T13:23:29.056'
$ws.Cells.Item(2, 4).Formula = '="281474991205262"'
$ws.Cells.Item(2, 4).Copy() | Out-Null
$ws.Cells.Item(2, 4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(2, 5).Formula = '="132"'
$ws.Cells.Item(2, 5).Copy() | Out-Null
$ws.Cells.Item(2, 5).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(2, 6).Formula = '="52215867"'
$ws.Cells.Item(2, 6).Copy() | Out-Null
$ws.Cells.Item(2, 6).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(2, 7).Value = 'EMMANUEL SALCEDO'
$ws.Cells.Item(2, 8).Value = 20.65046053
$ws.Cells.Item(2, 9).Value = -103.35116976
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276206556/TJaAJ2AJpb-camera-video-segment-driver-1743276209056.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSIR3L2ROE%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T060003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBwaCXVzLXdlc3QtMiJHMEUCIFQxHFeXs7MSUF7Q1VN2ptHYXePVkMaXFbSNef84375UAiEAmnlEh3RaQHtc3JHswC19tB8f9d3UZhfMVgZk5mrWj1oq5gMIhf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDKdweOOfBVtKgKyFiyq6A%2FuKeG6ES6dtY%2Bwx8dZCwEY2UfKDI5Ec8nn5mWE12dn%2F4nF0tPBSUi4jHSkP7G1WMAZZwzZPQW4zqx11wv61VEPMB2fcZvkoz5XqvfLGWCsb%2B6LEXvLvyRrnRt8JZb7DUmTsbHT8oDnwROjLk%2BZGEqmxtRYYHfSSNTxu8q6oApNnkroyAz0zEVyLhZtTmQIRITy3xcCUxfONkcIa29eyWQ327PP0WfdXvb7cPFAiES5f1eysOJnP5f3NI0WskyeVghikyzvgyKprnwNBkZ3id2R2H4u4D1B4oZLtmuBZ3ZjiKoy6oLJm26gwnq92QZI%2FbTQFlRpZS%2Fx9lSEikZyuoTz0AFUNWLXumQVXKMI%2BqsnJo%2BtDznSpU1BpDpWkan5a6PAMLDECzyoAJf%2FjZdxYqaotSthrYP0cihMZ8yibqw6aeCECGmtk7AwT7PhN4xiCsxiFYtzHWgm%2BanDc3L8ueBGzjCt%2FSZDpz5TVaJjJBKX9deS%2BoW2rmv2nTwdHdOF2hiDXZiiA%2Fw5EvaJsrNVuFfthf32rWtBdY%2BW%2BqEx1cI3E9asqo9tX61E6B8oLPs4pJ1DaNvl8Gompz1kwyoujvwY6pQGLcjqDqeAGmwjypKc5ms4bmhOjZBml8gF6Uz1EG2xtcJ5RrB%2FsdgFi%2Bv6HeRtx%2FDs%2BfxdoQMAMCroMkuLNaCrIEL%2FzYtza5%2BicK%2BBveS8876ptcWyrgAVdN97SooJMyJ2a1xEzIXyU9nKSkb7LRUh3qLRzIQ7KDDpv4h0IQeZuyY9hxNi1DTQV5mOaYwWyifWCLV7LXNxwiU3jDXrZOLIp81G5Lk8%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2014%3A00%3A03%20GMT&X-Amz-Signature=9f4109884f7707823073e17af7b30eb6c4f5b6d66ed17cf13c4f1291dc0e0140'
$ws.Cells.Item(2, 12).Value = 'No video URL'

# Row 3
$ws.Cells.Item(3, 1).Value = '281474991205262-1743276167363'
$ws.Cells.Item(3, 2).Value = 'Mobile Usage'
$ws.Cells.Item(3, 3).Value = '2025-03-29T13:22:47.363'
$ws.Cells.Item(3, 4).Formula = '="281474991205262"'
$ws.Cells.Item(3, 4).Copy() | Out-Null
$ws.Cells.Item(3, 4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(3, 5).Formula = '="132"'
$ws.Cells.Item(3, 5).Copy() | Out-Null
$ws.Cells.Item(3, 5).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(3, 6).Formula = '="52215867"'
$ws.Cells.Item(3, 6).Copy() | Out-Null
$ws.Cells.Item(3, 6).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(3, 7).Value = 'EMMANUEL SALCEDO'
$ws.Cells.Item(3, 8).Value = 20.65263968
$ws.Cells.Item(3, 9).Value = -103.34923
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276164863/OF6jgx1Lcp-camera-video-segment-driver-1743276167363.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSIR3L2ROE%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T060003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBwaCXVzLXdlc3QtMiJHMEUCIFQxHFeXs7MSUF7Q1VN2ptHYXePVkMaXFbSNef84375UAiEAmnlEh3RaQHtc3JHswC19tB8f9d3UZhfMVgZk5mrWj1oq5gMIhf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDKdweOOfBVtKgKyFiyq6A%2FuKeG6ES6dtY%2Bwx8dZCwEY2UfKDI5Ec8nn5mWE12dn%2F4nF0tPBSUi4jHSkP7G1WMAZZwzZPQW4zqx11wv61VEPMB2fcZvkoz5XqvfLGWCsb%2B6LEXvLvyRrnRt8JZb7DUmTsbHT8oDnwROjLk%2BZGEqmxtRYYHfSSNTxu8q6oApNnkroyAz0zEVyLhZtTmQIRITy3xcCUxfONkcIa29eyWQ327PP0WfdXvb7cPFAiES5f1eysOJnP5f3NI0WskyeVghikyzvgyKprnwNBkZ3id2R2H4u4D1B4oZLtmuBZ3ZjiKoy6oLJm26gwnq92QZI%2FbTQFlRpZS%2Fx9lSEikZyuoTz0AFUNWLXumQVXKMI%2BqsnJo%2BtDznSpU1BpDpWkan5a6PAMLDECzyoAJf%2FjZdxYqaotSthrYP0cihMZ8yibqw6aeCECGmtk7AwT7PhN4xiCsxiFYtzHWgm%2BanDc3L8ueBGzjCt%2FSZDpz5TVaJjJBKX9deS%2BoW2rmv2nTwdHdOF2hiDXZiiA%2Fw5EvaJsrNVuFfthf32rWtBdY%2BW%2BqEx1cI3E9asqo9tX61E6B8oLPs4pJ1DaNvl8Gompz1kwyoujvwY6pQGLcjqDqeAGmwjypKc5ms4bmhOjZBml8gF6Uz1EG2xtcJ5RrB%2FsdgFi%2Bv6HeRtx%2FDs%2BfxdoQMAMCroMkuLNaCrIEL%2FzYtza5%2BicK%2BBveS8876ptcWyrgAVdN97SooJMyJ2a1xEzIXyU9nKSkb7LRUh3qLRzIQ7KDDpv4h0IQeZuyY9hxNi1DTQV5mOaYwWyifWCLV7LXNxwiU3jDXrZOLIp81G5Lk8%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2014%3A00%3A03%20GMT&X-Amz-Signature=7bff9165149817ce2d2364d8f496a9480a80250df5a051fc8cc0c225bb71f7e6'
$ws.Cells.Item(3, 12).Value = 'No video URL'

# Row 4
$ws.Cells.Item(4, 1).Value = '281474991205262-1743276111238'
$ws.Cells.Item(4, 2).Value = 'Mobile Usage'
$ws.Cells.Item(4, 3).Value = '2025-03-29T13:21:51.238'
$ws.Cells.Item(4, 4).Formula = '="281474991205262"'
$ws.Cells.Item(4, 4).Copy() | Out-Null
$ws.Cells.Item(4, 4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(4, 5).Formula = '="132"'
$ws.Cells.Item(4, 5).Copy() | Out-Null
$ws.Cells.Item(4, 5).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(4, 6).Formula = '="52215867"'
$ws.Cells.Item(4, 6).Copy() | Out-Null
$ws.Cells.Item(4, 6).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(4, 7).Value = 'EMMANUEL SALCEDO'
$ws.Cells.Item(4, 8).Value = 20.65489757
$ws.Cells.Item(4, 9).Value = -103.345146729
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276108738/gAxSKLP636-camera-video-segment-driver-1743276111238.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSIR3L2ROE%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T060003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBwaCXVzLXdlc3QtMiJHMEUCIFQxHFeXs7MSUF7Q1VN2ptHYXePVkMaXFbSNef84375UAiEAmnlEh3RaQHtc3JHswC19tB8f9d3UZhfMVgZk5mrWj1oq5gMIhf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDKdweOOfBVtKgKyFiyq6A%2FuKeG6ES6dtY%2Bwx8dZCwEY2UfKDI5Ec8nn5mWE12dn%2F4nF0tPBSUi4jHSkP7G1WMAZZwzZPQW4zqx11wv61VEPMB2fcZvkoz5XqvfLGWCsb%2B6LEXvLvyRrnRt8JZb7DUmTsbHT8oDnwROjLk%2BZGEqmxtRYYHfSSNTxu8q6oApNnkroyAz0zEVyLhZtTmQIRITy3xcCUxfONkcIa29eyWQ327PP0WfdXvb7cPFAiES5f1eysOJnP5f3NI0WskyeVghikyzvgyKprnwNBkZ3id2R2H4u4D1B4oZLtmuBZ3ZjiKoy6oLJm26gwnq92QZI%2FbTQFlRpZS%2Fx9lSEikZyuoTz0AFUNWLXumQVXKMI%2BqsnJo%2BtDznSpU1BpDpWkan5a6PAMLDECzyoAJf%2FjZdxYqaotSthrYP0cihMZ8yibqw6aeCECGmtk7AwT7PhN4xiCsxiFYtzHWgm%2BanDc3L8ueBGzjCt%2FSZDpz5TVaJjJBKX9deS%2BoW2rmv2nTwdHdOF2hiDXZiiA%2Fw5EvaJsrNVuFfthf32rWtBdY%2BW%2BqEx1cI3E9asqo9tX61E6B8oLPs4pJ1DaNvl8Gompz1kwyoujvwY6pQGLcjqDqeAGmwjypKc5ms4bmhOjZBml8gF6Uz1EG2xtcJ5RrB%2FsdgFi%2Bv6HeRtx%2FDs%2BfxdoQMAMCroMkuLNaCrIEL%2FzYtza5%2BicK%2BBveS8876ptcWyrgAVdN97SooJMyJ2a1xEzIXyU9nKSkb7LRUh3qLRzIQ7KDDpv4h0IQeZuyY9hxNi1DTQV5mOaYwWyifWCLV7LXNxwiU3jDXrZOLIp81G5Lk8%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2014%3A00%3A03%20GMT&X-Amz-Signature=28323b6feeddee7cd91fed749c451bedcd8668db02d7f68a51645038d0cd12c9'
$ws.Cells.Item(4, 12).Value = 'No video URL'

# Row 5
$ws.Cells.Item(5, 1).Value = '281474991205341-1743265655240'
$ws.Cells.Item(5, 2).Value = 'No Seat Belt'
$ws.Cells.Item(5, 3).Value = '2025-03-29T10:27:35.240'
$ws.Cells.Item(5, 4).Formula = '="281474991205341"'
$ws.Cells.Item(5, 4).Copy() | Out-Null
$ws.Cells.Item(5, 4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(5, 5).Formula = '="140"'
$ws.Cells.Item(5, 5).Copy() | Out-Null
$ws.Cells.Item(5, 5).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(5, 6).Formula = '="51834149"'
$ws.Cells.Item(5, 6).Copy() | Out-Null
$ws.Cells.Item(5, 6).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(5, 7).Value = 'ABRAHAM ARANA'
$ws.Cells.Item(5, 8).Value = 20.70084996
$ws.Cells.Item(5, 9).Value = -103.4157048
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205341/1743265652740/hq4516Mzws-camera-video-segment-driver-1743265655240.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSIR3L2ROE%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T060003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBwaCXVzLXdlc3QtMiJHMEUCIFQxHFeXs7MSUF7Q1VN2ptHYXePVkMaXFbSNef84375UAiEAmnlEh3RaQHtc3JHswC19tB8f9d3UZhfMVgZk5mrWj1oq5gMIhf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDKdweOOfBVtKgKyFiyq6A%2FuKeG6ES6dtY%2Bwx8dZCwEY2UfKDI5Ec8nn5mWE12dn%2F4nF0tPBSUi4jHSkP7G1WMAZZwzZPQW4zqx11wv61VEPMB2fcZvkoz5XqvfLGWCsb%2B6LEXvLvyRrnRt8JZb7DUmTsbHT8oDnwROjLk%2BZGEqmxtRYYHfSSNTxu8q6oApNnkroyAz0zEVyLhZtTmQIRITy3xcCUxfONkcIa29eyWQ327PP0WfdXvb7cPFAiES5f1eysOJnP5f3NI0WskyeVghikyzvgyKprnwNBkZ3id2R2H4u4D1B4oZLtmuBZ3ZjiKoy6oLJm26gwnq92QZI%2FbTQFlRpZS%2Fx9lSEikZyuoTz0AFUNWLXumQVXKMI%2BqsnJo%2BtDznSpU1BpDpWkan5a6PAMLDECzyoAJf%2FjZdxYqaotSthrYP0cihMZ8yibqw6aeCECGmtk7AwT7PhN4xiCsxiFYtzHWgm%2BanDc3L8ueBGzjCt%2FSZDpz5TVaJjJBKX9deS%2BoW2rmv2nTwdHdOF2hiDXZiiA%2Fw5EvaJsrNVuFfthf32rWtBdY%2BW%2BqEx1cI3E9asqo9tX61E6B8oLPs4pJ1DaNvl8Gompz1kwyoujvwY6pQGLcjqDqeAGmwjypKc5ms4bmhOjZBml8gF6Uz1EG2xtcJ5RrB%2FsdgFi%2Bv6HeRtx%2FDs%2BfxdoQMAMCroMkuLNaCrIEL%2FzYtza5%2BicK%2BBveS8876ptcWyrgAVdN97SooJMyJ2a1xEzIXyU9nKSkb7LRUh3qLRzIQ7KDDpv4h0IQeZuyY9hxNi1DTQV5mOaYwWyifWCLV7LXNxwiU3jDXrZOLIp81G5Lk8%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2014%3A00%3A03%20GMT&X-Amz-Signature=61dbd8b6e70cd24b9832d1690e95e49dda48c0ba13e983303e05c5ea8e1dc8c3'
$ws.Cells.Item(5, 12).Value = 'No video URL'

# Row 6
$ws.Cells.Item(6, 1).Value = '281474991205262-1743256002469'
$ws.Cells.Item(6, 2).Value = 'No Seat Belt'
$ws.Cells.Item(6, 3).Value = '2025-03-29T07:46:42.469'
$ws.Cells.Item(6, 4).Formula = '="281474991205262"'
$ws.Cells.Item(6, 4).Copy() | Out-Null
$ws.Cells.Item(6, 4).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(6, 5).Formula = '="132"'
$ws.Cells.Item(6, 5).Copy() | Out-Null
$ws.Cells.Item(6, 5).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(6, 6).Formula = '="52215867"'
$ws.Cells.Item(6, 6).Copy() | Out-Null
$ws.Cells.Item(6, 6).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(6, 7).Value = 'EMMANUEL SALCEDO'
$ws.Cells.Item(6, 8).Value = 20.65084239
$ws.Cells.Item(6, 9).Value = -103.35091326
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743255999969/Yjq7ZuVo4p-camera-video-segment-driver-1743256002469.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSIR3L2ROE%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T060003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBwaCXVzLXdlc3QtMiJHMEUCIFQxHFeXs7MSUF7Q1VN2ptHYXePVkMaXFbSNef84375UAiEAmnlEh3RaQHtc3JHswC19tB8f9d3UZhfMVgZk5mrWj1oq5gMIhf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDKdweOOfBVtKgKyFiyq6A%2FuKeG6ES6dtY%2Bwx8dZCwEY2UfKDI5Ec8nn5mWE12dn%2F4nF0tPBSUi4jHSkP7G1WMAZZwzZPQW4zqx11wv61VEPMB2fcZvkoz5XqvfLGWCsb%2B6LEXvLvyRrnRt8JZb7DUmTsbHT8oDnwROjLk%2BZGEqmxtRYYHfSSNTxu8q6oApNnkroyAz0zEVyLhZtTmQIRITy3xcCUxfONkcIa29eyWQ327PP0WfdXvb7cPFAiES5f1eysOJnP5f3NI0WskyeVghikyzvgyKprnwNBkZ3id2R2H4u4D1B4oZLtmuBZ3ZjiKoy6oLJm26gwnq92QZI%2FbTQFlRpZS%2Fx9lSEikZyuoTz0AFUNWLXumQVXKMI%2BqsnJo%2BtDznSpU1BpDpWkan5a6PAMLDECzyoAJf%2FjZdxYqaotSthrYP0cihMZ8yibqw6aeCECGmtk7AwT7PhN4xiCsxiFYtzHWgm%2BanDc3L8ueBGzjCt%2FSZDpz5TVaJjJBKX9deS%2BoW2rmv2nTwdHdOF2hiDXZiiA%2Fw5EvaJsrNVuFfthf32rWtBdY%2BW%2BqEx1cI3E9asqo9tX61E6B8oLPs4pJ1DaNvl8Gompz1kwyoujvwY6pQGLcjqDqeAGmwjypKc5ms4bmhOjZBml8gF6Uz1EG2xtcJ5RrB%2FsdgFi%2Bv6HeRtx%2FDs%2BfxdoQMAMCroMkuLNaCrIEL%2FzYtza5%2BicK%2BBveS8876ptcWyrgAVdN97SooJMyJ2a1xEzIXyU9nKSkb7LRUh3qLRzIQ7KDDpv4h0IQeZuyY9hxNi1DTQV5mOaYwWyifWCLV7LXNxwiU3jDXrZOLIp81G5Lk8%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2014%3A00%3A03%20GMT&X-Amz-Signature=88c25607a944708932062d5f1c9024156c456255814375390de1451a6d80b915'
$ws.Cells.Item(6, 12).Value = 'No video URL'

$excel.CutCopyMode = $false
